$d = $word.ActiveDocument

function Insert-ParaAfter {
    # Inserts a new paragraph right after $d.Paragraphs.Item($anchorIndex)
    # and puts $text into it (unless $text is empty, in which case the
    # paragraph is left empty). Returns the index of the new paragraph.
    param($anchorIndex, [string]$text)
    $p = $d.Paragraphs.Item($anchorIndex)
    $p.Range.InsertParagraphAfter()
    $newIndex = $anchorIndex + 1
    if ($text -ne "") {
        $np = $d.Paragraphs.Item($newIndex)
        $np.Range.Text = $text
    }
    return $newIndex
}

function Replace-InRange {
    param($range, [string]$old, [string]$new)
    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# We work from the bottom of the document upward so paragraph indices
# for content we have not processed yet stay valid.

# --- Paragraph 16: "IPF = (... его релиза) / LOC" ----------------------
$p16 = $d.Paragraphs.Item(16)
Replace-InRange $p16.Range `
    "IPF = (Число обнаруженных ошибок до выпуска его релиза) / LOС" `
    "IPF = (Число обнаруженных ошибок до выпуска релиза) / строки в коде"

$idx = 16
$idx = Insert-ParaAfter $idx ""
$idx = Insert-ParaAfter $idx "Минимальное значение – 1. Чем больше – тем лучше"
$idx = Insert-ParaAfter $idx "Значение де-факто – 0,05"
$idx = Insert-ParaAfter $idx ""

# --- Paragraph 13: "T = Общее время ..." (text unchanged) --------------
$idx = 13
$idx = Insert-ParaAfter $idx "Минимально допустимое значение – 1. Чем меньше – тем лучше"
$idx = Insert-ParaAfter $idx "Значение де-факто – 1,5"

# --- Paragraph 10: empty paragraph before "Качества ..." heading -------
$p10 = $d.Paragraphs.Item(10)
$p10.Range.InsertBefore("Максимально допустимое значение – 1. Чем меньше – тем лучше")
$idx = Insert-ParaAfter 10 "Значение де-факто – 0,5(30 секунд)"

# --- Paragraph 9: "PRR = Количество минут на перенесение карт ..." -----
$p9 = $d.Paragraphs.Item(9)
Replace-InRange $p9.Range `
    "PRR = Количество минут на перенесение карт в базу данных;" `
    "PRR = Количество минут на перенесение одной карты в базу данных;"

# --- Paragraph 6: "FS = (... )" (text unchanged) ------------------------
$idx = Insert-ParaAfter 6 "Максимальное значение – 1. Чем больше – тем лучше"
$idx = Insert-ParaAfter $idx "Значение де-факто – 0,9"

# --- Paragraph 3: "Productivity = LOC / ..." -----------------------------
$p3 = $d.Paragraphs.Item(3)
Replace-InRange $p3.Range `
    " = LOC / Рабочее время, затраченное на проект;" `
    " = строки в коде / Рабочее время, затраченное на проект"

$idx = 3
$idx = Insert-ParaAfter $idx "Максимальное значение – 100. Чем больше – тем лучше"
$idx = Insert-ParaAfter $idx "Значение де-факто -- 166  (на проект было отведено 30 часов, строки 10 000)"

Write-Host "Final paragraph count:" $d.Paragraphs.Count
